# Update the "想去人数" (F column) values on both the "展览" and "全部类型"
# worksheets, which carry identical data tables.

$wb = $excel.ActiveWorkbook

# Map of row number -> new value for column F
$updates = @{
    3  = 2975
    7  = 1624
    10 = 28
    11 = 1335
    13 = 456
    14 = 340
    17 = 119
    18 = 87
    19 = 101
    20 = 3064
    21 = 372
    22 = 89
    24 = 85
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
